# CORE_holdings.xlsx update
# - Roll the "as of" date in the confidentiality notice from 2021-03-19 to 2021-03-22
# - Refresh the Weight (D) / Percent Change (E) figures for rows 2-8
#
# The sheet ships protected, so we have to unprotect it before editing the
# locked cells and re-protect it again afterwards to leave it in the same
# (protected) state it was found in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Confidentiality notice: bump the "as of" date
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-22 for illustrative purposes only and are subject to change."

# Row 2 (IVE)
$ws.Range("D2").Value = 0.4971145069486708
$ws.Range("E2").Value = -0.00007113387395085713

# Row 3 (IVW)
$ws.Range("D3").Value = 0.2405365813572818
$ws.Range("E3").Value = 0.01460427135678399

# Row 4 (IJK)
$ws.Range("D4").Value = 0.09866813042593348
$ws.Range("E4").Value = 0.001411161000641314

# Row 5 (IJJ)
$ws.Range("D5").Value = 0.1034831456384883
$ws.Range("E5").Value = -0.01194058829239875

# Row 6 (IJS)
$ws.Range("D6").Value = 0.03134029409265042
$ws.Range("E6").Value = -0.0203865756322722

# Row 7 (IJT)
$ws.Range("D7").Value = 0.02885734153697536
$ws.Range("E7").Value = -0.005973807153251243

# Row 8 (Total) - only the Percent Change changes, Weight stays 1
$ws.Range("E8").Value = 0.001569777335935552

$ws.Protect()
